$wb = $excel.ActiveWorkbook

# Sheet "展览" (worksheet 1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 373
$ws1.Range("F3").Value = 806
$ws1.Range("F4").Value = 277
$ws1.Range("F5").Value = 936
$ws1.Range("F6").Value = 2248
$ws1.Range("F7").Value = 198

# Sheet "全部类型" (worksheet 4): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 373
$ws4.Range("F3").Value = 806
$ws4.Range("F4").Value = 277
$ws4.Range("F7").Value = 936
$ws4.Range("F8").Value = 2248
$ws4.Range("F10").Value = 198
